$d = $word.ActiveDocument

# Replace the first occurrence of $old (found via Find) with $new, preserving the
# run's formatting and avoiding the smart-quote autocorrect that Find.Execute's
# replace argument triggers (we assign Range.Text directly instead).
function Set-RangeText($old, $new) {
    $rng = $d.Content
    $found = $rng.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Text not found: $old"
    }
    $rng.Text = $new
}

# Locate $anchor (via Find), collapse to its end and return that collapsed Range,
# so callers can chain further InsertAfter calls onto it in document order.
function Get-EndRange($anchor) {
    $rng = $d.Content
    $found = $rng.Find.Execute($anchor, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Anchor not found: $anchor"
    }
    $rng.Collapse(0)
    return $rng
}

# =========================================================================
# Title
# =========================================================================
Set-RangeText "Unraveling the Quantum Landscape" "The Ethereal Realm of Chemistry: A Journey into the Intricacies of Matter"

# =========================================================================
# Author name
# =========================================================================
Set-RangeText " Ethan Winters" " Eleanor Thompson"

# =========================================================================
# Email address (split across runs: "winterse@quantumphysics" + "." + "edu")
# becomes "eleanor" + "." + "thompson@scholarlyedu" + new "." + new "org"
# =========================================================================
Set-RangeText "winterse@quantumphysics" "eleanor"
Set-RangeText "edu" "thompson@scholarlyedu"
$r = Get-EndRange "thompson@scholarlyedu"
$r.InsertAfter(".")
$r.Collapse(0)
$r.InsertAfter("org")

# =========================================================================
# Body paragraph 1 (first block of text, before the two <w:br/> line breaks)
# =========================================================================
Set-RangeText "The realm of quantum physics, a perplexing universe of interconnectedness and uncertainty, continues to captivate scientists and philosophers alike" "In the vast tapestry of science, chemistry stands as a beacon of discovery, illuminating the enigmatic world of matter"
Set-RangeText " As the foundation of modern physics, quantum mechanics has illuminated the infinitesimal world's intricate workings, revealing a realm where particles dance in probabilistic waves and probabilities govern the fabric of reality" " As a high school teacher, I aim to unveil the secrets of this captivating subject, guiding my students on an enthralling journey through the intricate dance of elements and molecules"
Set-RangeText " Delving into the quantum landscape signifies a journey into the profound depths of matter, energy, and the fundamental laws that weave the tapestry of existence" " Chemistry is a symphony of reactions and transformations, a saga of particles interacting in a harmonious ballet of creation and destruction"

# New sentences appended in order, each its own run followed by its own "." run
$r = Get-EndRange " Chemistry is a symphony of reactions and transformations, a saga of particles interacting in a harmonious ballet of creation and destruction"
$r.InsertAfter(".")
$r.Collapse(0)
$r.InsertAfter(" Each element, with its unique properties, weaves its magic, forming the very fabric of our universe")
$r.Collapse(0)
$r.InsertAfter(".")
$r.Collapse(0)
$r.InsertAfter(" It is a realm where the mysteries of the material world unravel, revealing the fundamental forces that shape our lives")
$r.Collapse(0)
$r.InsertAfter(".")
$r.Collapse(0)
$r.InsertAfter(" It is a subject that combines intellectual rigor with endless fascination, a testament to the boundless curiosity of the human spirit")

# =========================================================================
# Body paragraph 1, second block (after the first "<w:br/><w:br/>")
# =========================================================================
Set-RangeText "In this multifaceted domain, particles exhibit both wave-like and particle-like characteristics, a duality that defies classical intuition" "This captivating field unveils the secrets of how matter is composed, how it changes, and how it interacts with its surroundings"
Set-RangeText " The enigmatic nature of quantum entanglement further astounds, as particles separated by vast distances remain bound in an inseparable embrace of shared fate" " Unraveling the intricacies of these interactions, we glimpse the profound elegance of the natural world"
Set-RangeText " As we delve deeper into this microscopic realm, the uncertainty principle unveils the inherent interconnectedness of measurements, revealing the inherent limitations of our knowledge" " Chemistry is a gateway to understanding the very essence of things, from the smallest atoms to the grandest molecules"
Set-RangeText " The quantum world is a tapestry of phenomena that challenge our most fundamental assumptions about reality, inviting us to rethink the very nature of space, time, and existence itself" " As we delve deeper into this realm of discovery, we learn to harness the power of chemical reactions to create new materials, devise novel medicines, and address some of the world's most pressing challenges"

# =========================================================================
# Body paragraph 1, third block (after the second "<w:br/><w:br/>")
# =========================================================================
Set-RangeText "The study of quantum physics is an endeavor punctuated by both awe and perplexity" "The study of chemistry is not merely an academic pursuit; it is an endeavor that connects us to the world around us"
Set-RangeText " Its profound implications have rippled across numerous fields, from computation and communication to cosmology and biology" " By understanding the fundamental principles that govern chemical processes, we gain insights into a myriad of phenomena, from the vibrant colors of flowers to the intricate workings of our own bodies"
Set-RangeText " As scientists continue to unravel the enigmatic tapestry of quantum mechanics, we glimpse the potential for transformative technologies, insights into the nature of consciousness, and a deeper understanding of the universe's fundamental nature" " Chemistry empowers us to make informed decisions about our health, our environment, and our future"

# New sentences appended in order, each its own run followed by its own "." run
$r = Get-EndRange " Chemistry empowers us to make informed decisions about our health, our environment, and our future"
$r.InsertAfter(".")
$r.Collapse(0)
$r.InsertAfter(" It is a discipline that fosters critical thinking, problem-solving skills, and a deep appreciation for the natural world")
$r.Collapse(0)
$r.InsertAfter(".")
$r.Collapse(0)
$r.InsertAfter(" It prepares us to navigate an increasingly complex world where scientific literacy is essential for informed citizenship")

# =========================================================================
# Summary heading text is unchanged ("Summary"); the <w:lastRenderedPageBreak/>
# that now precedes it is a pagination artifact recalculated automatically.
# =========================================================================

# =========================================================================
# Summary paragraph
# =========================================================================
Set-RangeText "The realm of quantum physics presents a mind-boggling universe of interconnectedness and uncertainty, challenging classical notions of reality" "Chemistry, an alluring field of scientific exploration, unveils the mysteries of matter, its composition, and its interactions"
Set-RangeText " Particles exhibit wave-like and particle-like behaviors, entanglements defy distance, and the uncertainty principle reveals the interconnectedness of measurements" " Through the study of chemistry, we gain profound insights into the fundamental forces that shape our universe and the intricate workings of the natural world"

# The old two runs (" Quantum physics has profound implications ... the " and
# "cosmos and consciousness") collapse into a single new sentence.
Set-RangeText " Quantum physics has profound implications across diverse fields, inspiring new technologies and reshaping our understanding of the " " It empowers us with knowledge and skills essential for navigating an increasingly complex world, enabling us to address global challenges and make informed decisions about our health, our environment, and our future"
Set-RangeText "cosmos and consciousness" ""

Set-RangeText " As scientists continue to unravel the quantum landscape, the possibilities for transformative advancements are boundless" " Chemistry is a gateway to understanding the very essence of things and fosters a deep appreciation for the harmonious ballet of elements and molecules that comprise our existence"
